$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells F1:H1, matching the style used by existing headers (A1:E1)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

$ws.Range("A1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null

# New boolean columns F, G, H for rows 2-12 (all FALSE except G9 = TRUE)
$ws.Range("F2:H12").Value = $false
$ws.Range("G9").Value = $true
